# The "Min. ord.kvant." column (column H) was removed from the sheet.
# Deleting the entire column shifts the following columns (I,J,K,L) one
# position to the left (I->H, J->I, K->J, L->K), taking their values,
# formulas and styles with them, and removes the now-unused
# "Min. ord.kvant." entry from the shared strings table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(8).EntireColumn.Delete()
